$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values are stored as text (not numbers) in the source workbook, e.g.
# "1,725,882.00" must remain literal text, not be parsed into a numeric value.
# Setting NumberFormat to Text ("@") before assigning the value keeps these
# cells as plain text without leaving a "quote prefix" style behind on save.
$row2 = $ws.Range("A2:E2")
$row2.NumberFormat = "@"

$ws.Range("A2").Value = "1,725,882.00"
$ws.Range("B2").Value = "383,493.00"
$ws.Range("C2").Value = "1,200,000.00"
$ws.Range("D2").Value = "909,375.00"
$ws.Range("E2").Value = "383,493.00"

$row2.Style = "Normal"
